$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.146.64"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = "'2.264.17"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'153.27"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +15,207.29%  '
$ws.Range('D6').Value = "'306.01"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('D7').Value = "'94.44"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.60%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('D11').Value = "'33.25"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.09%  '
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').Value = "'6.66"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('D15').Value = "'2.616.58"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').Value = "'14.35"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.53%  '
$ws.Range('D17').Value = "'2.264.89"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = "'0.786"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.64%  '
$ws.Range('D19').Value = "'41.992.69"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('D20').Value = "'12.65"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.49%  '
$ws.Range('D21').Value = "'0.0₃0915"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.08%  '
$ws.Range('D22').Value = "'5.99"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('D23').Value = "'68.14"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('D24').Value = "'243.87"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('E25').Value = '  +1.87%  '
$ws.Range('E26').Value = '  +2.19%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').Value = "'23.99"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').Value = "'9.69"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').Value = "'34.98"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.50%  '
$ws.Range('D32').Value = "'159.93"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('E33').Value = '  +3.53%  '
$ws.Range('D34').Value = "'0.999"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = "'0.0745"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').Value = "'3.07"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').Value = "'17.03"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.56%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').Value = "'2.37"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = "'0.105"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').Value = "'1.80"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('D42').Value = "'4.09"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.35%  '
$ws.Range('D43').Value = "'2.002.84"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.05%  '
$ws.Range('D44').Value = "'19.52"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('E45').Value = '  +10.70%  '
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('D47').Value = "'10.20"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').Value = "'2.91"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('D49').Value = "'53.65"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.77%  '
$ws.Range('D50').Value = "'72.89"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.92%  '
$ws.Range('E51').Value = '  +0.10%  '
